$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily data rows (2021-08-24 .. 2021-09-01), appended after the last
# existing row (357) — columns: A=date serial, B=nuovi pos., C=somma mobile
# 7gg., D=somma mobile 7gg. per 100mila abitanti.
$newRows = @(
    @{ Row = 358; Date = 44432; B = 0;  C = 23; D = 69.73289270230119 },
    @{ Row = 359; Date = 44433; B = 1;  C = 23; D = 69.73289270230119 },
    @{ Row = 360; Date = 44434; B = 1;  C = 24; D = 72.76475760240123 },
    @{ Row = 361; Date = 44435; B = 12; C = 26; D = 78.82848740260134 },
    @{ Row = 362; Date = 44436; B = 1;  C = 24; D = 72.76475760240123 },
    @{ Row = 363; Date = 44437; B = 13; C = 35; D = 106.1152715035018 },
    @{ Row = 364; Date = 44438; B = 10; C = 38; D = 115.210866203802  },
    @{ Row = 365; Date = 44439; B = 0;  C = 38; D = 115.210866203802  },
    @{ Row = 366; Date = 44440; B = 0;  C = 37; D = 112.1790013037019 }
)

foreach ($r in $newRows) {
    $rowNum = $r.Row
    $prevRow = $rowNum - 1

    # Carry over column A's date style (border/bold/center/date numfmt) from
    # the row directly above instead of constructing a brand new style.
    $ws.Range("A$prevRow").Copy() | Out-Null
    $ws.Range("A$rowNum").PasteSpecial(-4122) | Out-Null

    $ws.Range("A$rowNum").Value = $r.Date
    $ws.Range("B$rowNum").Value = $r.B
    $ws.Range("C$rowNum").Value = $r.C
    $ws.Range("D$rowNum").Value = $r.D
}
